$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New stopword abbreviations appended to the list. The shared-string table
# records new unique strings in the order they are first entered, so "yg"
# is entered before "dgn" (even though "dgn" lands in the earlier row, A94)
# to reproduce the exact shared-string ordering of the target workbook.
$ws.Range("A95").Value = "yg"
$ws.Range("A94").Value = "dgn"
$ws.Range("A96").Value = "hny"

# Move the selection to the next empty cell below the newly entered data,
# matching the cursor position left behind after typing the new rows.
$ws.Range("A97").Select()
